# AZ2522021|3:11PM Adding Parameters for All navigation Urls
#
# The "Reports" row label in the SimpleSearch parameters sheet had an
# accidental double underscore ("Reports__Tab_URL"). Fix the label text and
# leave the user's selection on that cell (matching the saved sheetView
# selection in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SimpleSearch")
$ws.Activate()

# Correct the typo'd navigation-url parameter name in column A (row 10).
$ws.Range("A10").Value = "Reports_Tab_URL"

# Reflect the author's final selection/cursor position on the sheet.
$ws.Range("A10").Select()
